$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.131.12"
$ws.Range("E2").Value = "  +5.08%  "
$ws.Range("D3").Value = "2.588.03"
$ws.Range("E3").Value = "  +6.64%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "505.82"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").Value = "155.80"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "0.585"
$ws.Range("E8").Value = "  -5.13%  "
$ws.Range("D9").Value = "2.620.54"
$ws.Range("E9").Value = "  +7.07%  "
$ws.Range("D10").Value = "6.44"
$ws.Range("E10").Value = "  +4.67%  "
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +3.96%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "3.072.80"
$ws.Range("E14").Value = "  +8.04%  "
$ws.Range("D15").Value = "60.289.00"
$ws.Range("E15").Value = "  +5.31%  "
$ws.Range("D16").Value = "21.62"
$ws.Range("E16").Value = "  +4.97%  "
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +4.58%  "
$ws.Range("D18").Value = "2.626.24"
$ws.Range("E18").Value = "  +7.46%  "
$ws.Range("D19").Value = "4.76"
$ws.Range("E19").Value = "  +2.90%  "
$ws.Range("D20").Value = "343.59"
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("D21").Value = "10.38"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "6.14"
$ws.Range("E22").Value = "  +3.91%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "59.99"
$ws.Range("E24").Value = "  +3.73%  "
$ws.Range("D25").Value = "0.422"
$ws.Range("E25").Value = "  +5.08%  "
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "0.0₃0850"
$ws.Range("E28").Value = "  +8.22%  "
$ws.Range("D29").Value = "7.53"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "156.85"
$ws.Range("E31").Value = "  +3.92%  "
$ws.Range("D32").Value = "19.31"
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("E33").Value = "  +3.11%  "
$ws.Range("D34").Value = "5.70"
$ws.Range("E34").Value = "  +7.52%  "
$ws.Range("D35").Value = "4.00"
$ws.Range("E35").Value = "  +5.55%  "
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("D37").Value = "306.72"
$ws.Range("E37").Value = "  +7.90%  "
$ws.Range("D38").Value = "0.845"
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  +7.02%  "
$ws.Range("D40").Value = "3.76"
$ws.Range("E40").Value = "  +6.65%  "
$ws.Range("E41").Value = "  +26.91%  "
$ws.Range("E42").Value = "  +4.52%  "
$ws.Range("D43").Value = "0.625"
$ws.Range("E43").Value = "  +3.98%  "
$ws.Range("D44").Value = "0.0570"
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("D45").Value = "0.100"
$ws.Range("E45").Value = "  -1.29%  "
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("D47").Value = "19.84"
$ws.Range("E47").Value = "  +12.75%  "
$ws.Range("D48").Value = "4.87"
$ws.Range("E48").Value = "  +6.91%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +3.33%  "
$ws.Range("D50").Value = "2.042.70"
$ws.Range("E50").Value = "  +7.57%  "
$ws.Range("D51").Value = "10.27"
$ws.Range("E51").Value = "  +0.62%  "

Write-Host "Updated cryptos list"
